# Apply the "add w12 - node, python" edit to the "Parcial 2" sheet.
#
# Summary of the change:
#  - Insert two new columns before the old column P (new week-12 columns),
#    which pushes every column from old P..AB to new R..AD.
#  - Resize a handful of columns (D, E, F, I, L) that were manually
#    resized by the author alongside the insert.
#  - Fill in the new week-12 header cells (P1/Q1) with "node-redis" /
#    "python-redis".
#  - Fill in the new week-12 data column (now P/Q) with 0s for every
#    student row.
#  - Mark week 8-11 (L:O) as done (1) for row 2 (Alarico Mercado Vázquez)
#    and mark week 9-10 (M:N) as done for row 4, whose repo link (I4) is
#    updated to the real GitHub URL.
#  - Update the selection / active cell shown when the sheet was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parcial 2")
$ws.Activate()

# --- 1. Insert two new columns before the old "P" column -------------------
$ws.Range("P1:Q1").EntireColumn.Insert() | Out-Null

# --- 2. Column width tweaks made alongside the insert -----------------------
$ws.Range("D1").EntireColumn.ColumnWidth = 4.33203125
$ws.Range("E1").EntireColumn.ColumnWidth = 1.6640625
$ws.Range("F1").EntireColumn.ColumnWidth = 2.83203125
$ws.Range("I1").EntireColumn.ColumnWidth = 26
$ws.Range("L1").EntireColumn.ColumnWidth = 10.6640625

# --- 3. New week-12 header text (shared strings "node-redis"/"python-redis") ---
$ws.Range("P1").Value = "node-redis"
$ws.Range("Q1").Value = "python-redis"

# --- 4. New week-12 data column: default to 0 for every student row (2-23) ---
$ws.Range("P2:Q23").Value = 0

# --- 5. Row 2 (Alarico Mercado Vázquez): mark weeks 8-11 as complete --------
$ws.Range("L2:O2").Value = 1

# --- 6. Row 4: update repo link and mark weeks 9-10 as complete ------------
$ws.Range("I4").Value = "https://github.com/blackhawk42/administracion"
$ws.Range("M4:N4").Value = 1

# --- 7. Restore the active cell selection used when saving -----------------
$ws.Range("M1").Select() | Out-Null
